{"js": "const pairs = [\n  [\"945\u00f78=\", \"913\u00f73=\"],\n  [\"844\u00f74=\", \"376\u00f72=\"],\n  [\"705\u00f79=\", \"557\u00f76=\"],\n  [\"954\u00f75=\", \"951\u00f75=\"],\n  [\"278\u00f78=\", \"941\u00f79=\"],\n  [\"252\u00f78=\", \"681\u00f78=\"],\n  [\"366\u00f74=\", \"231\u00f76=\"],\n  [\"322\u00f75=\", \"159\u00f78=\"],\n  [\"934\u00f78=\", \"211\u00f77=\"],\n  [\"800\u00f79=\", \"591\u00f77=\"],\n  [\"556\u00f72=\", \"829\u00f72=\"],\n  [\"164\u00f74=\", \"345\u00f76=\"],\n  [\"594\u00f79=\", \"930\u00f72=\"],\n  [\"748\u00f74=\", \"835\u00f77=\"],\n  [\"317\u00f75=\", \"849\u00f75=\"],\n  [\"666\u00f75=\", \"956\u00f72=\"],\n  [\"413\u00f75=\", \"170\u00f73=\"],\n  [\"961\u00f77=\", \"124\u00f73=\"],\n  [\"994\u00f72=\", \"209\u00f78=\"],\n  [\"253\u00f73=\", \"716\u00f72=\"],\n  [\"463\u00f79=\", \"814\u00f77=\"],\n  [\"470\u00f72=\", \"959\u00f72=\"],\n  [\"426\u00f73=\", \"365\u00f77=\"],\n  [\"373\u00f73=\", \"805\u00f79=\"],\n  [\"847\u00f76=\", \"510\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ Old = \"945\u00f78=\"; New = \"913\u00f73=\" },\n    @{ Old = \"844\u00f74=\"; New = \"376\u00f72=\" },\n    @{ Old = \"705\u00f79=\"; New = \"557\u00f76=\" },\n    @{ Old = \"954\u00f75=\"; New = \"951\u00f75=\" },\n    @{ Old = \"278\u00f78=\"; New = \"941\u00f79=\" },\n    @{ Old = \"252\u00f78=\"; New = \"681\u00f78=\" },\n    @{ Old = \"366\u00f74=\"; New = \"231\u00f76=\" },\n    @{ Old = \"322\u00f75=\"; New = \"159\u00f78=\" },\n    @{ Old = \"934\u00f78=\"; New = \"211\u00f77=\" },\n    @{ Old = \"800\u00f79=\"; New = \"591\u00f77=\" },\n    @{ Old = \"556\u00f72=\"; New = \"829\u00f72=\" },\n    @{ Old = \"164\u00f74=\"; New = \"345\u00f76=\" },\n    @{ Old = \"594\u00f79=\"; New = \"930\u00f72=\" },\n    @{ Old = \"748\u00f74=\"; New = \"835\u00f77=\" },\n    @{ Old = \"317\u00f75=\"; New = \"849\u00f75=\" },\n    @{ Old = \"666\u00f75=\"; New = \"956\u00f72=\" },\n    @{ Old = \"413\u00f75=\"; New = \"170\u00f73=\" },\n    @{ Old = \"961\u00f77=\"; New = \"124\u00f73=\" },\n    @{ Old = \"994\u00f72=\"; New = \"209\u00f78=\" },\n    @{ Old = \"253\u00f73=\"; New = \"716\u00f72=\" },\n    @{ Old = \"463\u00f79=\"; New = \"814\u00f77=\" },\n    @{ Old = \"470\u00f72=\"; New = \"959\u00f72=\" },\n    @{ Old = \"426\u00f73=\"; New = \"365\u00f77=\" },\n    @{ Old = \"373\u00f73=\"; New = \"805\u00f79=\" },\n    @{ Old = \"847\u00f76=\"; New = \"510\u00f73=\" }\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Execute(\n        $pair.Old,   # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap = wdFindContinue\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace = wdReplaceAll\n    )\n}\n"}
